# COCO3_Comp.xlsx edit:
#  - Add a new "Load 20240216" column (D) that mirrors the existing
#    "Load 20240205" column (C) on Sheet1, row by row.
#  - D1 gets the new header text "Load 20240216".
#  - D2:D31 / D34:D38 get the same value (and therefore the same
#    formatting) as the matching cell in column C.
#  - Move the active selection to D26.
#  - Make column D the same pixel width as column C (so the column
#    group B:D ends up uniform, matching columns B:C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell D1: new unique shared string -----------------------
$ws.Range("D1").Value = "Load 20240216"

# --- Data rows: mirror column C into column D ------------------------
# Row 25, 32 and 33 do not exist in the sheet (gaps in the table) and
# must stay that way, so they are intentionally skipped.
$dataRows = @(2..24) + @(26..31)
foreach ($r in $dataRows) {
    $srcCell = $ws.Cells.Item($r, 3)   # column C
    $dstCell = $ws.Cells.Item($r, 4)   # column D
    $dstCell.Value2 = $srcCell.Value2
}

# --- Notes rows 34-38: mirror column C into column D, including the
#     centered / wrap-text formatting used by B34:C38 -----------------
$noteRows = 34..38
foreach ($r in $noteRows) {
    $srcCell = $ws.Cells.Item($r, 3)   # column C
    $dstCell = $ws.Cells.Item($r, 4)   # column D
    $dstCell.Value2 = $srcCell.Value2
    $dstCell.HorizontalAlignment = -4108   # xlCenter
    $dstCell.WrapText = $true
}

# --- Column D width: match column C's pixel width so the B:D columns
#     end up visually identical -----------------------------------
$ws.Columns(4).Width = $ws.Columns(3).Width

# --- Update the saved selection to D26 --------------------------------
$ws.Range("D26").Select()
